$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 10506
$ws.Range("J69").Value = 13999
$ws.Range("L69").Value = 41997
$ws.Range("N69").Value = -43745
$ws.Range("H70").Value = 4292.857
$ws.Range("I70").Value = 4067.3333
$ws.Range("K70").Value = 12201.9999
$ws.Range("M70").Value = -11931.9999
$ws.Range("H72").Value = 10506
$ws.Range("J72").Value = 13999
$ws.Range("L72").Value = 125991
$ws.Range("N72").Value = -134727
$ws.Range("H73").Value = 4292.857
$ws.Range("I73").Value = 4067.3333
$ws.Range("K73").Value = 12201.9999
$ws.Range("M73").Value = -11265.9999
$ws.Range("H88").Value = 3669.25
$ws.Range("I88").Value = 6000
$ws.Range("J88").Value = 2892.3333
$ws.Range("K88").Value = 6000
$ws.Range("L88").Value = 2892.3333
$ws.Range("M88").Value = -5594
$ws.Range("N88").Value = -3704.3333
$ws.Range("H91").Value = 3669.25
$ws.Range("I91").Value = 6000
$ws.Range("J91").Value = 2892.3333
$ws.Range("K91").Value = 6000
$ws.Range("L91").Value = 2892.3333
$ws.Range("M91").Value = -4596
$ws.Range("N91").Value = -5700.3333
$ws.Range("H113").Value = 2701.6667
$ws.Range("I113").Value = 1802.5
$ws.Range("K113").Value = 1802.5
$ws.Range("M113").Value = 1451.5
$ws.Range("H138").Value = 3034.7334
$ws.Range("I138").Value = 2021.55
$ws.Range("J138").Value = 3845.28
$ws.Range("K138").Value = 6064.65
$ws.Range("L138").Value = 11535.84
$ws.Range("M138").Value = -924.6499999999996
$ws.Range("N138").Value = -21815.84

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5409870
$ws.Range("I32").Value = 6670976.5
$ws.Range("K32").Value = 6670976.5
$ws.Range("M32").Value = -6670689.5
$ws.Range("H61").Value = 4712.2666
$ws.Range("I61").Value = 3379.6155
$ws.Range("K61").Value = 3379.6155
$ws.Range("M61").Value = -3167.6155
$ws.Range("H74").Value = 2797.8708
$ws.Range("I74").Value = 1280.5834
$ws.Range("K74").Value = 1280.5834
$ws.Range("M74").Value = -406.5834
$ws.Range("H77").Value = 2797.8708
$ws.Range("I77").Value = 1280.5834
$ws.Range("K77").Value = 6402.916999999999
$ws.Range("M77").Value = -2034.916999999999
$ws.Range("H122").Value = 2516.3333
$ws.Range("I122").Value = 2206
$ws.Range("J122").Value = 4999
$ws.Range("K122").Value = 6618
$ws.Range("L122").Value = 14997
$ws.Range("M122").Value = -4168
$ws.Range("N122").Value = -19897
$ws.Range("H132").Value = 6709.25
$ws.Range("I132").Value = 4057.5715
$ws.Range("J132").Value = 14664.286
$ws.Range("K132").Value = 12172.7145
$ws.Range("L132").Value = 43992.858
$ws.Range("M132").Value = -9642.7145
$ws.Range("N132").Value = -49052.858
$ws.Range("H136").Value = 4712.2666
$ws.Range("I136").Value = 3379.6155
$ws.Range("K136").Value = 10138.8465
$ws.Range("M136").Value = -7588.8465

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2490.3845
$ws.Range("I86").Value = 2281.25
$ws.Range("J86").Value = 5000
$ws.Range("K86").Value = 2281.25
$ws.Range("L86").Value = 5000
$ws.Range("M86").Value = -1158.25
$ws.Range("N86").Value = -7246
$ws.Range("H89").Value = 2490.3845
$ws.Range("I89").Value = 2281.25
$ws.Range("J89").Value = 5000
$ws.Range("K89").Value = 11406.25
$ws.Range("L89").Value = 25000
$ws.Range("M89").Value = -5790.25
$ws.Range("N89").Value = -36232
$ws.Range("H99").Value = 951
$ws.Range("I99").Value = 951
$ws.Range("K99").Value = 951
$ws.Range("M99").Value = 547
$ws.Range("H107").Value = 2923.0667
$ws.Range("I107").Value = 1674.6923
$ws.Range("K107").Value = 1674.6923
$ws.Range("M107").Value = 245.3077000000001
$ws.Range("H134").Value = 5126.1665
$ws.Range("I134").Value = 2265.92
$ws.Range("K134").Value = 6797.76
$ws.Range("M134").Value = -4262.76

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5154.0264
$ws.Range("I31").Value = 3146.923
$ws.Range("J31").Value = 6197.72
$ws.Range("K31").Value = 3146.923
$ws.Range("L31").Value = 6197.72
$ws.Range("M31").Value = -2851.923
$ws.Range("N31").Value = -6787.72
$ws.Range("H34").Value = 5154.0264
$ws.Range("I34").Value = 3146.923
$ws.Range("J34").Value = 6197.72
$ws.Range("K34").Value = 3146.923
$ws.Range("L34").Value = 6197.72
$ws.Range("M34").Value = -2944.923
$ws.Range("N34").Value = -6601.72
$ws.Range("H55").Value = 24872.25
$ws.Range("J55").Value = 24872.25
$ws.Range("L55").Value = 24872.25
$ws.Range("N55").Value = -25502.25
$ws.Range("H58").Value = 5638.3237
$ws.Range("I58").Value = 4530
$ws.Range("K58").Value = 4530
$ws.Range("M58").Value = -4327
$ws.Range("H100").Value = 74125
$ws.Range("J100").Value = 74125
$ws.Range("L100").Value = 74125
$ws.Range("N100").Value = -76289
$ws.Range("H111").Value = 60000
$ws.Range("J111").Value = 60000
$ws.Range("L111").Value = 60000
$ws.Range("N111").Value = -68180
$ws.Range("H136").Value = 5638.3237
$ws.Range("I136").Value = 4530
$ws.Range("K136").Value = 13590
$ws.Range("M136").Value = -11040

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 5980.75
$ws.Range("J2").Value = 8881.75
$ws.Range("L2").Value = 53290.5
$ws.Range("N2").Value = -53516.5
$ws.Range("H8").Value = 355
$ws.Range("I8").Value = 355
$ws.Range("K8").Value = 1065
$ws.Range("M8").Value = -926
$ws.Range("H9").Value = 1764.5714
$ws.Range("I9").Value = 639
$ws.Range("K9").Value = 1917
$ws.Range("M9").Value = -1693
$ws.Range("H19").Value = 361
$ws.Range("I19").Value = 383
$ws.Range("J19").Value = 295
$ws.Range("K19").Value = 1149
$ws.Range("L19").Value = 885
$ws.Range("M19").Value = -975
$ws.Range("N19").Value = -1233
$ws.Range("H54").Value = 290
$ws.Range("I54").Value = 290
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 870
$ws.Range("L54").Value = 0
$ws.Range("M54").Value = -311
$ws.Range("N54").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H40").Value = 25999
$ws.Range("J40").Value = 25999
$ws.Range("L40").Value = 25999
$ws.Range("N40").Value = -26301
$ws.Range("H102").Value = 1029.8125
$ws.Range("I102").Value = 816.0909
$ws.Range("K102").Value = 816.0909
$ws.Range("M102").Value = 805.9091
$ws.Range("H123").Value = 42954.617
$ws.Range("J123").Value = 44267.777
$ws.Range("L123").Value = 44267.777
$ws.Range("N123").Value = -49167.777
$ws.Range("H132").Value = 15088.363
$ws.Range("J132").Value = 15599
$ws.Range("L132").Value = 46797
$ws.Range("N132").Value = -51857

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("I25").Value = 11644471
$ws.Range("J25").Value = 250065000
$ws.Range("K25").Value = 11644471
$ws.Range("L25").Value = 250065000
$ws.Range("M25").Value = -11644241
$ws.Range("N25").Value = -250065460
$ws.Range("H40").Value = 33338428
$ws.Range("I40").Value = 38466496
$ws.Range("K40").Value = 38466496
$ws.Range("M40").Value = -38466360
$ws.Range("H68").Value = 2788.0908
$ws.Range("I68").Value = 2752.8572
$ws.Range("J68").Value = 2849.75
$ws.Range("K68").Value = 2752.8572
$ws.Range("L68").Value = 2849.75
$ws.Range("M68").Value = -2003.8572
$ws.Range("N68").Value = -4347.75
$ws.Range("H71").Value = 2788.0908
$ws.Range("I71").Value = 2752.8572
$ws.Range("J71").Value = 2849.75
$ws.Range("K71").Value = 13764.286
$ws.Range("L71").Value = 14248.75
$ws.Range("M71").Value = -10020.286
$ws.Range("N71").Value = -21736.75
$ws.Range("H100").Value = 8337012
$ws.Range("I100").Value = 11907851
$ws.Range("K100").Value = 11907851
$ws.Range("M100").Value = -11907310

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H40").Value = 49997.5
$ws.Range("I40").Value = 49997.5
$ws.Range("K40").Value = 49997.5
$ws.Range("M40").Value = -49848.5
$ws.Range("H43").Value = 67506.75
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").ClearContents()
$ws.Range("H62").Value = 27285.715
$ws.Range("I62").Value = 48666.668
$ws.Range("J62").Value = 11250
$ws.Range("K62").Value = 48666.668
$ws.Range("L62").Value = 11250
$ws.Range("M62").Value = -48042.668
$ws.Range("N62").Value = -12498
$ws.Range("H65").Value = 27285.715
$ws.Range("I65").Value = 48666.668
$ws.Range("J65").Value = 11250
$ws.Range("K65").Value = 243333.34
$ws.Range("L65").Value = 56250
$ws.Range("M65").Value = -240213.34
$ws.Range("N65").Value = -62490
$ws.Range("H81").Value = 2064.4688
$ws.Range("I81").Value = 1339.5
$ws.Range("J81").Value = 7139.25
$ws.Range("K81").Value = 2679
$ws.Range("L81").Value = 14278.5
$ws.Range("M81").Value = -1618
$ws.Range("N81").Value = -16400.5
$ws.Range("H84").Value = 2064.4688
$ws.Range("I84").Value = 1339.5
$ws.Range("J84").Value = 7139.25
$ws.Range("K84").Value = 13395
$ws.Range("L84").Value = 71392.5
$ws.Range("M84").Value = -8091
$ws.Range("N84").Value = -82000.5
